$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("rights_and_functions")

# Insert a new row above row 45, shifting rows 45+ down by one
$ws.Rows("45").Insert()

# Populate new row 45
$ws.Cells.Item(45, 2).Formula = '="micration_x_to_"&TEXT(TODAY(),"jjjj")&"_"&TEXT(TODAY(),"MM")&"_"&TEXT(TODAY(),"TT")&".sql"'
$ws.Cells.Item(45, 3).Value = "template_micration.sql"

# Give the (now empty) cell two rows below a date number format
$ws.Cells.Item(47, 2).NumberFormat = "mm-dd-yy"

$ws.Range("B47").Select()
